$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old leading "序" / numbering column (column A) entirely,
# shifting every other column one to the left (B->A, C->B, ... G->F).
$ws.Range("A1").EntireColumn.Delete()

# Strip the redundant "（白带）" suffix from the item-name column (now column B).
$ws.Range("B2").Value = "PH"
$ws.Range("B3").Value = "清洁度"
$ws.Range("B4").Value = "上皮细胞"
$ws.Range("B5").Value = "乳酸杆菌"
$ws.Range("B6").Value = "白细胞"
$ws.Range("B7").Value = "红细胞"
$ws.Range("B8").Value = "线索细胞"
$ws.Range("B9").Value = "阴道毛滴虫"
$ws.Range("B10").Value = "霉菌"
$ws.Range("B11").Value = "疑似支原体"
$ws.Range("B12").Value = "疑似衣原体"

# Normalise the double-dash ranges to a single dash (now column F).
$ws.Range("F2").Value = "4-4.5"
$ws.Range("F6").Value = "0-15"

# Clean up the mangled code values in column A for the last two rows.
$ws.Range("A11").Value = "bdysyyt"
$ws.Range("A12").Value = "bdysyyt"
